# Apply the cryptos-list update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain-text numbers (e.g. "29.915.33",
# "1.000", "84.88"). Force the column to Text format before writing so
# Excel does not auto-coerce values like "1.000" into the number 1.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.915.33"
$ws.Range("E2").Value = "  +6.28%  "
# Row 3
$ws.Range("D3").Value = "1.877.12"
# Row 4
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.08%  "
# Row 5
$ws.Range("D5").Value = "248.45"
$ws.Range("E5").Value = "  +2.02%  "
# Row 6
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.05%  "
# Row 7
$ws.Range("D7").Value = "0.4985"
$ws.Range("E7").Value = "  +1.31%  "
# Row 8
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").Value = "45.87"
$ws.Range("E8").Value = "  +8.74%  "
# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.2852"
$ws.Range("E9").Value = "  +6.81%  "
# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.06537"
$ws.Range("E10").Value = "  +4.45%  "
# Row 11
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.880.48"
$ws.Range("E11").Value = "  +5.70%  "
# Row 12
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "17.05"
$ws.Range("E12").Value = "  +3.66%  "
# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07230"
$ws.Range("E13").Value = "  +2.80%  "
# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.6667"
$ws.Range("E14").Value = "  +6.24%  "
# Row 15
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "84.88"
$ws.Range("E15").Value = "  +6.33%  "
# Row 16
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "4.815"
$ws.Range("E16").Value = "  +3.73%  "
# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "29.923.50"
$ws.Range("E17").Value = "  +6.37%  "
# Row 18
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.03%  "
# Row 19
$ws.Range("D19").Value = "12.86"
$ws.Range("E19").Value = "  +6.86%  "
# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.000007505"
$ws.Range("E20").Value = "  +3.66%  "
# Row 21
$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  -0.11%  "
# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.122.43"
$ws.Range("E22").Value = "  +5.64%  "
# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "4.753"
$ws.Range("E23").Value = "  +4.17%  "
# Row 24
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "5.528"
$ws.Range("E24").Value = "  +5.76%  "
# Row 25
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "9.005"
$ws.Range("E25").Value = "  +3.18%  "
# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "145.09"
$ws.Range("E26").Value = "  +2.71%  "
# Row 27
$ws.Range("B27").Value = "BitcoinCash"
$ws.Range("C27").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D27").Value = "134.90"
$ws.Range("E27").Value = "  +23.37%  "
# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "16.68"
$ws.Range("E28").Value = "  +5.87%  "
# Row 29
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "1.949"
$ws.Range("E29").Value = "  +4.90%  "
# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "1.372"
$ws.Range("E30").Value = "  -0.97%  "
# Row 31
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "4.188"
$ws.Range("E31").Value = "  +0.29%  "
# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.08620"
$ws.Range("E32").Value = "  +4.40%  "
# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "3.875"
$ws.Range("E33").Value = "  +2.98%  "
# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.05074"
$ws.Range("E34").Value = "  +3.95%  "
# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.128"
$ws.Range("E35").Value = "  +5.43%  "
# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6843"
$ws.Range("E36").Value = "  +5.04%  "
# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.697"
$ws.Range("E37").Value = "  +3.16%  "
# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "2.286"
$ws.Range("E38").Value = "  +11.92%  "
# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.739"
$ws.Range("E39").Value = "  +5.67%  "
# Row 40
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "0.9560"
$ws.Range("E40").Value = "  +1.03%  "
# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.01632"
$ws.Range("E41").Value = "  +5.42%  "
# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "6.098"
$ws.Range("E42").Value = "  +3.52%  "
# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "104.48"
$ws.Range("E43").Value = "  +4.69%  "
# Row 44
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.01%  "
# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.4202"
$ws.Range("E45").Value = "  +5.60%  "
# Row 46
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "7.438"
$ws.Range("E46").Value = "  +3.72%  "
# Row 47
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1252"
$ws.Range("E47").Value = "  +3.50%  "
# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05628"
$ws.Range("E48").Value = "  +3.57%  "
# Row 49
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "32.36"
$ws.Range("E49").Value = "  +5.28%  "
# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "8.216"
$ws.Range("E50").Value = "  +2.85%  "
# Row 51
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.3700"
$ws.Range("E51").Value = "  +6.53%  "

